# The GSC export rolls its 89-day window forward by one day:
#   - the oldest date (2025-10-26) and its "Valid" count drop off the top
#   - every remaining date/value shifts up one row
#   - a new date (2026-01-23) with its "Valid" count is appended at the bottom
#
# Deleting row 2 shifts rows 3:90 up to 2:89 (dates + values move together),
# and the engine re-derives the shared-string table automatically (dropping
# the now-unused "2025-10-26" entry). We then append the new trailing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day - shifts all rows below it up by one.
$ws.Rows.Item(2).Delete()

# Append the new day at the bottom (row 90 after the shift).
# Force text formatting first so the yyyy-MM-dd-looking string isn't
# auto-converted into a date serial number.
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = "2026-01-23"
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 24

# Re-apply the plain (unstyled) format used by the rest of the date column,
# so the new cell matches its neighbours instead of keeping the "@" style.
$ws.Range("A89").Copy()
$ws.Range("A90").PasteSpecial(-4122)
